$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Cells.Item(116, 8).Value = 13467.223
$ws.Cells.Item(116, 9).Value = 13467.223
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 13467.223
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -10025.223
$ws.Cells.Item(116, 14).Value = $null
# Row 129
$ws.Cells.Item(129, 8).Value = 23810334
$ws.Cells.Item(129, 9).Value = 111111520
$ws.Cells.Item(129, 10).Value = 919.57574
$ws.Cells.Item(129, 11).Value = 333334560
$ws.Cells.Item(129, 12).Value = 2758.72722
$ws.Cells.Item(129, 13).Value = -333329560
$ws.Cells.Item(129, 14).Value = -12758.72722
# Row 132
$ws.Cells.Item(132, 8).Value = 2264.0312
$ws.Cells.Item(132, 9).Value = 2295.5
$ws.Cells.Item(132, 10).Value = 2127.6667
$ws.Cells.Item(132, 11).Value = 6886.5
$ws.Cells.Item(132, 12).Value = 6383.000100000001
$ws.Cells.Item(132, 13).Value = -4356.5
$ws.Cells.Item(132, 14).Value = -11443.0001
# Row 137
$ws.Cells.Item(137, 8).Value = 2439.85
$ws.Cells.Item(137, 9).Value = 1686.2667
$ws.Cells.Item(137, 10).Value = 4700.6
$ws.Cells.Item(137, 11).Value = 5058.800099999999
$ws.Cells.Item(137, 12).Value = 14101.8
$ws.Cells.Item(137, 13).Value = -2508.800099999999
$ws.Cells.Item(137, 14).Value = -19201.8

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5312.646
$ws.Cells.Item(32, 9).Value = 3465.279
$ws.Cells.Item(32, 10).Value = 21200
$ws.Cells.Item(32, 11).Value = 3465.279
$ws.Cells.Item(32, 12).Value = 21200
$ws.Cells.Item(32, 13).Value = -3178.279
$ws.Cells.Item(32, 14).Value = -21774
# Row 61
$ws.Cells.Item(61, 8).Value = 5376.1
$ws.Cells.Item(61, 9).Value = 4293.2915
$ws.Cells.Item(61, 10).Value = 9707.333000000001
$ws.Cells.Item(61, 11).Value = 4293.2915
$ws.Cells.Item(61, 12).Value = 9707.333000000001
$ws.Cells.Item(61, 13).Value = -4081.2915
$ws.Cells.Item(61, 14).Value = -10131.333
# Row 74
$ws.Cells.Item(74, 8).Value = 2297.8096
$ws.Cells.Item(74, 9).Value = 2374.2856
$ws.Cells.Item(74, 10).Value = 2144.8572
$ws.Cells.Item(74, 11).Value = 2374.2856
$ws.Cells.Item(74, 12).Value = 2144.8572
$ws.Cells.Item(74, 13).Value = -1500.2856
$ws.Cells.Item(74, 14).Value = -3892.8572
# Row 77
$ws.Cells.Item(77, 8).Value = 2297.8096
$ws.Cells.Item(77, 9).Value = 2374.2856
$ws.Cells.Item(77, 10).Value = 2144.8572
$ws.Cells.Item(77, 11).Value = 11871.428
$ws.Cells.Item(77, 12).Value = 10724.286
$ws.Cells.Item(77, 13).Value = -7503.428
$ws.Cells.Item(77, 14).Value = -19460.286
# Row 132
$ws.Cells.Item(132, 8).Value = 2591.3784
$ws.Cells.Item(132, 9).Value = 1182.05
$ws.Cells.Item(132, 10).Value = 4249.4116
$ws.Cells.Item(132, 11).Value = 3546.15
$ws.Cells.Item(132, 12).Value = 12748.2348
$ws.Cells.Item(132, 13).Value = -1016.15
$ws.Cells.Item(132, 14).Value = -17808.2348
# Row 136
$ws.Cells.Item(136, 8).Value = 5376.1
$ws.Cells.Item(136, 9).Value = 4293.2915
$ws.Cells.Item(136, 10).Value = 9707.333000000001
$ws.Cells.Item(136, 11).Value = 12879.8745
$ws.Cells.Item(136, 12).Value = 29121.999
$ws.Cells.Item(136, 13).Value = -10329.8745
$ws.Cells.Item(136, 14).Value = -34221.999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = $null
$ws.Cells.Item(132, 14).Value = $null
# Row 134
$ws.Cells.Item(134, 8).Value = 2482.3333
$ws.Cells.Item(134, 9).Value = 1778.826
$ws.Cells.Item(134, 10).Value = 4100.4
$ws.Cells.Item(134, 11).Value = 5336.478
$ws.Cells.Item(134, 12).Value = 12301.2
$ws.Cells.Item(134, 13).Value = -2801.478
$ws.Cells.Item(134, 14).Value = -17371.2

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 10340.028
$ws.Cells.Item(31, 9).Value = 1204.3158
$ws.Cells.Item(31, 10).Value = 21188.688
$ws.Cells.Item(31, 11).Value = 1204.3158
$ws.Cells.Item(31, 12).Value = 21188.688
$ws.Cells.Item(31, 13).Value = -909.3158000000001
$ws.Cells.Item(31, 14).Value = -21778.688
# Row 32
$ws.Cells.Item(32, 8).Value = 5340
$ws.Cells.Item(32, 9).Value = 5340
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 5340
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = -5024
$ws.Cells.Item(32, 14).Value = $null
# Row 34
$ws.Cells.Item(34, 8).Value = 10340.028
$ws.Cells.Item(34, 9).Value = 1204.3158
$ws.Cells.Item(34, 10).Value = 21188.688
$ws.Cells.Item(34, 11).Value = 1204.3158
$ws.Cells.Item(34, 12).Value = 21188.688
$ws.Cells.Item(34, 13).Value = -1002.3158
$ws.Cells.Item(34, 14).Value = -21592.688
# Row 58
$ws.Cells.Item(58, 8).Value = 1934.8667
$ws.Cells.Item(58, 9).Value = 763.25
$ws.Cells.Item(58, 10).Value = 2360.9092
$ws.Cells.Item(58, 11).Value = 763.25
$ws.Cells.Item(58, 12).Value = 2360.9092
$ws.Cells.Item(58, 13).Value = -560.25
$ws.Cells.Item(58, 14).Value = -2766.9092
# Row 99
$ws.Cells.Item(99, 8).Value = 2007352.5
$ws.Cells.Item(99, 9).Value = 3201166
$ws.Cells.Item(99, 10).Value = 17663.334
$ws.Cells.Item(99, 11).Value = 3201166
$ws.Cells.Item(99, 12).Value = 17663.334
$ws.Cells.Item(99, 13).Value = -3199668
$ws.Cells.Item(99, 14).Value = -20659.334
# Row 126
$ws.Cells.Item(126, 8).Value = 2007352.5
$ws.Cells.Item(126, 9).Value = 3201166
$ws.Cells.Item(126, 10).Value = 17663.334
$ws.Cells.Item(126, 11).Value = 9603498
$ws.Cells.Item(126, 12).Value = 52990.00199999999
$ws.Cells.Item(126, 13).Value = -9601028
$ws.Cells.Item(126, 14).Value = -57930.00199999999
# Row 132
$ws.Cells.Item(132, 8).Value = 2634.6
$ws.Cells.Item(132, 9).Value = 1855
$ws.Cells.Item(132, 10).Value = 5753
$ws.Cells.Item(132, 11).Value = 5565
$ws.Cells.Item(132, 12).Value = 17259
$ws.Cells.Item(132, 13).Value = -3035
$ws.Cells.Item(132, 14).Value = -22319
# Row 134
$ws.Cells.Item(134, 8).Value = 9370.1
$ws.Cells.Item(134, 9).Value = 8960.875
$ws.Cells.Item(134, 10).Value = 11007
$ws.Cells.Item(134, 11).Value = 26882.625
$ws.Cells.Item(134, 12).Value = 33021
$ws.Cells.Item(134, 13).Value = -24347.625
$ws.Cells.Item(134, 14).Value = -38091
# Row 136
$ws.Cells.Item(136, 8).Value = 1934.8667
$ws.Cells.Item(136, 9).Value = 763.25
$ws.Cells.Item(136, 10).Value = 2360.9092
$ws.Cells.Item(136, 11).Value = 2289.75
$ws.Cells.Item(136, 12).Value = 7082.7276
$ws.Cells.Item(136, 13).Value = 260.25
$ws.Cells.Item(136, 14).Value = -12182.7276

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Cells.Item(4, 8).Value = 1270.7858
$ws.Cells.Item(4, 9).Value = 333.5
$ws.Cells.Item(4, 10).Value = 1973.75
$ws.Cells.Item(4, 11).Value = 1000.5
$ws.Cells.Item(4, 12).Value = 5921.25
$ws.Cells.Item(4, 13).Value = -888.5
$ws.Cells.Item(4, 14).Value = -6145.25
# Row 5
$ws.Cells.Item(5, 8).Value = 1169.5
$ws.Cells.Item(5, 9).Value = 531.8570999999999
$ws.Cells.Item(5, 10).Value = 2285.375
$ws.Cells.Item(5, 11).Value = 1595.5713
$ws.Cells.Item(5, 12).Value = 6856.125
$ws.Cells.Item(5, 13).Value = -1483.5713
$ws.Cells.Item(5, 14).Value = -7080.125
# Row 122
$ws.Cells.Item(122, 8).Value = 827.625
$ws.Cells.Item(122, 9).Value = 697.5454999999999
$ws.Cells.Item(122, 10).Value = 1113.8
$ws.Cells.Item(122, 11).Value = 6277.9095
$ws.Cells.Item(122, 12).Value = 10024.2
$ws.Cells.Item(122, 13).Value = -3827.9095
$ws.Cells.Item(122, 14).Value = -14924.2
# Row 132
$ws.Cells.Item(132, 8).Value = 997.13635
$ws.Cells.Item(132, 9).Value = 824.5
$ws.Cells.Item(132, 10).Value = 1457.5
$ws.Cells.Item(132, 11).Value = 7420.5
$ws.Cells.Item(132, 12).Value = 13117.5
$ws.Cells.Item(132, 13).Value = -4890.5
$ws.Cells.Item(132, 14).Value = -18177.5
# Row 135
$ws.Cells.Item(135, 8).Value = 1169.5
$ws.Cells.Item(135, 9).Value = 531.8570999999999
$ws.Cells.Item(135, 10).Value = 2285.375
$ws.Cells.Item(135, 11).Value = 4786.7139
$ws.Cells.Item(135, 12).Value = 20568.375
$ws.Cells.Item(135, 13).Value = -2251.7139
$ws.Cells.Item(135, 14).Value = -25638.375

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 1489
$ws.Cells.Item(102, 9).Value = 1302
$ws.Cells.Item(102, 10).Value = 1676
$ws.Cells.Item(102, 11).Value = 1302
$ws.Cells.Item(102, 12).Value = 1676
$ws.Cells.Item(102, 13).Value = 320
$ws.Cells.Item(102, 14).Value = -4920
# Row 113
$ws.Cells.Item(113, 8).Value = 1931.6666
$ws.Cells.Item(113, 9).Value = 1553.1428
$ws.Cells.Item(113, 10).Value = 3256.5
$ws.Cells.Item(113, 11).Value = 1553.1428
$ws.Cells.Item(113, 12).Value = 3256.5
$ws.Cells.Item(113, 13).Value = 616.8571999999999
$ws.Cells.Item(113, 14).Value = -7596.5
# Row 126
$ws.Cells.Item(126, 8).Value = 2058.7407
$ws.Cells.Item(126, 9).Value = 1623.6666
$ws.Cells.Item(126, 10).Value = 2602.5833
$ws.Cells.Item(126, 11).Value = 4870.9998
$ws.Cells.Item(126, 12).Value = 7807.749899999999
$ws.Cells.Item(126, 13).Value = -2400.9998
$ws.Cells.Item(126, 14).Value = -12747.7499

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 26187.045
$ws.Cells.Item(40, 9).Value = 30451.055
$ws.Cells.Item(40, 10).Value = 6999
$ws.Cells.Item(40, 11).Value = 30451.055
$ws.Cells.Item(40, 12).Value = 6999
$ws.Cells.Item(40, 13).Value = -30315.055
$ws.Cells.Item(40, 14).Value = -7271
# Row 132
$ws.Cells.Item(132, 8).Value = 6659.0356
$ws.Cells.Item(132, 9).Value = 9052.4375
$ws.Cells.Item(132, 10).Value = 3467.8333
$ws.Cells.Item(132, 11).Value = 27157.3125
$ws.Cells.Item(132, 12).Value = 10403.4999
$ws.Cells.Item(132, 13).Value = -24627.3125
$ws.Cells.Item(132, 14).Value = -15463.4999
# Row 136
$ws.Cells.Item(136, 8).Value = 2883.8293
$ws.Cells.Item(136, 9).Value = 2271.6765
$ws.Cells.Item(136, 10).Value = 5857.143
$ws.Cells.Item(136, 11).Value = 6815.029500000001
$ws.Cells.Item(136, 12).Value = 17571.429
$ws.Cells.Item(136, 13).Value = -4265.029500000001
$ws.Cells.Item(136, 14).Value = -22671.429

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 2523.7058
$ws.Cells.Item(132, 9).Value = 1821.125
$ws.Cells.Item(132, 10).Value = 3148.2222
$ws.Cells.Item(132, 11).Value = 5463.375
$ws.Cells.Item(132, 12).Value = 9444.6666
$ws.Cells.Item(132, 13).Value = -2933.375
$ws.Cells.Item(132, 14).Value = -14504.6666
# Row 136
$ws.Cells.Item(136, 8).Value = 22366.057
$ws.Cells.Item(136, 9).Value = 56126.723
$ws.Cells.Item(136, 10).Value = 5003.4287
$ws.Cells.Item(136, 11).Value = 168380.169
$ws.Cells.Item(136, 12).Value = 15010.2861
$ws.Cells.Item(136, 13).Value = -165830.169
$ws.Cells.Item(136, 14).Value = -20110.2861
